# Refresh market-board price snapshots (currentAveragePrice*, LevePrice*, LeveProfit*)
# for the Sargatanas Profits leve-crafting sheets. Values sourced from the latest
# Universalis price pull picked up by the scheduled runner.
$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 2: Mercury Rising / Quicksilver
$ws.Range("H2").Value = 42.608696
$ws.Range("I2").Value = 42.608696
$ws.Range("K2").Value = 42.608696
$ws.Range("M2").Value = 70.39130399999999
# Row 17: One for the Road / Potion
$ws.Range("H17").Value = 1390.72
$ws.Range("J17").Value = 1390.72
$ws.Range("L17").Value = 4172.16
$ws.Range("N17").Value = -4508.16
# Row 18: You Grow, Girl / Growth Formula Beta
$ws.Range("H18").Value = 1662.8
$ws.Range("I18").Value = 1768.75
$ws.Range("K18").Value = 1768.75
$ws.Range("M18").Value = -1484.75
# Row 28: The Writing Is Not on the Wall / Enchanted Silver Ink
$ws.Range("H28").Value = 1653.7333
$ws.Range("I28").Value = 1129.5
$ws.Range("J28").Value = 2702.2
$ws.Range("K28").Value = 1129.5
$ws.Range("L28").Value = 2702.2
$ws.Range("M28").Value = -644.5
$ws.Range("N28").Value = -3672.2
# Row 86: Filling in the Blanks / Enchanted Aurum Regis Ink
$ws.Range("H86").Value = 69446090
$ws.Range("I86").Value = 160714980
$ws.Range("J86").Value = 5557870
$ws.Range("K86").Value = 160714980
$ws.Range("L86").Value = 5557870
$ws.Range("M86").Value = -160713857
$ws.Range("N86").Value = -5560116
# Row 89: Ink into Antiquity (L) / Enchanted Aurum Regis Ink
$ws.Range("H89").Value = 69446090
$ws.Range("I89").Value = 160714980
$ws.Range("J89").Value = 5557870
$ws.Range("K89").Value = 803574900
$ws.Range("L89").Value = 27789350
$ws.Range("M89").Value = -803569284
$ws.Range("N89").Value = -27800582
# Row 92: Whinier than the Sword / Enchanted Koppranickel Ink
$ws.Range("H92").Value = 550.78125
$ws.Range("I92").Value = 446.26923
$ws.Range("J92").Value = 1003.6667
$ws.Range("K92").Value = 446.26923
$ws.Range("L92").Value = 1003.6667
$ws.Range("M92").Value = 801.73077
$ws.Range("N92").Value = -3499.6667
# Row 106: Making Your Mark / Enchanted Palladium Ink
$ws.Range("H106").Value = 2547
$ws.Range("I106").Value = 2637.5
$ws.Range("K106").Value = 2637.5
$ws.Range("M106").Value = -2006.5
# Row 113: Amaro Kart / Starch Glue
$ws.Range("H113").Value = 27791702
$ws.Range("I113").Value = 2484
$ws.Range("J113").Value = 38479864
$ws.Range("K113").Value = 2484
$ws.Range("L113").Value = 38479864
$ws.Range("M113").Value = 770
$ws.Range("N113").Value = -38486372
# Row 116: Growing Up / Growth Formula Kappa
$ws.Range("H116").Value = 12506127
$ws.Range("I116").Value = 27781234
$ws.Range("J116").Value = 8311.546
$ws.Range("K116").Value = 27781234
$ws.Range("L116").Value = 8311.546
$ws.Range("M116").Value = -27777792
$ws.Range("N116").Value = -15195.546
# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 1393267.4
$ws.Range("J138").Value = 1965905.6
$ws.Range("L138").Value = 5897716.800000001
$ws.Range("N138").Value = -5907996.800000001

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots / Bronze Ingot
$ws.Range("H2").Value = 1626.2858
$ws.Range("I2").Value = 1381.1052
$ws.Range("J2").Value = 2143.889
$ws.Range("K2").Value = 1381.1052
$ws.Range("L2").Value = 2143.889
$ws.Range("M2").Value = -1268.1052
$ws.Range("N2").Value = -2369.889
# Row 45: Hollow Hallmarks / Mythril Ingot
$ws.Range("H45").Value = 3997.2683
$ws.Range("I45").Value = 2867.3684
$ws.Range("J45").Value = 4973.091
$ws.Range("K45").Value = 2867.3684
$ws.Range("L45").Value = 4973.091
$ws.Range("M45").Value = -2490.3684
$ws.Range("N45").Value = -5727.091
# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 142863410
$ws.Range("J61").Value = 166672820
$ws.Range("L61").Value = 166672820
$ws.Range("N61").Value = -166673244
# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 24031.87
$ws.Range("I74").Value = 34836.066
$ws.Range("K74").Value = 34836.066
$ws.Range("M74").Value = -33962.066
# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 24031.87
$ws.Range("I77").Value = 34836.066
$ws.Range("K77").Value = 174180.33
$ws.Range("M77").Value = -169812.33
# Row 116: No Scope / Titanbronze Ingot
$ws.Range("H116").Value = 1626.2858
$ws.Range("I116").Value = 1381.1052
$ws.Range("J116").Value = 2143.889
$ws.Range("K116").Value = 1381.1052
$ws.Range("L116").Value = 2143.889
$ws.Range("M116").Value = 912.8948
$ws.Range("N116").Value = -6731.889
# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 142863410
$ws.Range("J136").Value = 166672820
$ws.Range("L136").Value = 500018460
$ws.Range("N136").Value = -500023560

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells / Bronze Ingot
$ws.Range("H3").Value = 1626.2858
$ws.Range("I3").Value = 1381.1052
$ws.Range("J3").Value = 2143.889
$ws.Range("K3").Value = 1381.1052
$ws.Range("L3").Value = 2143.889
$ws.Range("M3").Value = -1267.1052
$ws.Range("N3").Value = -2371.889
# Row 40: Can You Spare a Dolabra / Steel Dolabra
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
# Row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 500500.5
$ws.Range("I86").Value = 1000001
$ws.Range("J86").Value = 1000
$ws.Range("K86").Value = 1000001
$ws.Range("L86").Value = 1000
$ws.Range("M86").Value = -998878
$ws.Range("N86").Value = -3246
# Row 87: Winter Weather Conditions / Adamantite Dolabra
$ws.Range("H87").Value = 44166.668
$ws.Range("J87").Value = 44166.668
$ws.Range("L87").Value = 44166.668
$ws.Range("N87").Value = -46662.668
# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 500500.5
$ws.Range("I89").Value = 1000001
$ws.Range("J89").Value = 1000
$ws.Range("K89").Value = 5000005
$ws.Range("L89").Value = 5000
$ws.Range("M89").Value = -4994389
$ws.Range("N89").Value = -16232
# Row 90: The Nightsoil Is Dark and Full of Terrors (L) / Adamantite Dolabra
$ws.Range("H90").Value = 44166.668
$ws.Range("J90").Value = 44166.668
$ws.Range("L90").Value = 132500.004
$ws.Range("N90").Value = -144980.004
# Row 96: Hammer Time / High Steel Sledgehammer
$ws.Range("H96").Value = 14674.6
$ws.Range("I96").Value = 4857.75
$ws.Range("K96").Value = 4857.75
$ws.Range("M96").Value = -2111.75

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 99: O Pine / Pine Lumber
$ws.Range("H99").Value = 4161.1333
$ws.Range("I99").Value = 1689.875
$ws.Range("K99").Value = 1689.875
$ws.Range("M99").Value = -191.875
# Row 105: Zelkova, My Love / Zelkova Lumber
$ws.Range("H105").Value = 2646619.2
$ws.Range("I105").Value = 3402129.5
$ws.Range("K105").Value = 3402129.5
$ws.Range("M105").Value = -3400382.5
# Row 110: A Stronger Offense / Applewood Spear
$ws.Range("H110").Value = 59500
$ws.Range("J110").Value = 59500
$ws.Range("L110").Value = 59500
$ws.Range("N110").Value = -67680
# Row 126: A Better Conductor / Red Pine Lumber
$ws.Range("H126").Value = 4161.1333
$ws.Range("I126").Value = 1689.875
$ws.Range("K126").Value = 5069.625
$ws.Range("M126").Value = -2599.625
# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 4852.857
$ws.Range("I134").Value = 2853.9
$ws.Range("K134").Value = 8561.700000000001
$ws.Range("M134").Value = -6026.700000000001

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 2: Pork Is a Salty Food / Table Salt
$ws.Range("H2").Value = 100310.76
$ws.Range("I2").Value = 18871.875
$ws.Range("K2").Value = 113231.25
$ws.Range("M2").Value = -113118.25
# Row 37: I Love Lamprey / Eel Pie
$ws.Range("H37").Value = 76999
$ws.Range("J37").Value = 76999
$ws.Range("L37").Value = 230997
$ws.Range("N37").Value = -231221
# Row 97: The Frier Never Lies / Cottonseed Oil
$ws.Range("H97").Value = 473.85715
$ws.Range("J97").Value = 343
$ws.Range("L97").Value = 1029
$ws.Range("N97").Value = -2021
# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 1959.8334
$ws.Range("I131").Value = 1419.5
$ws.Range("J131").Value = 2230
$ws.Range("K131").Value = 4258.5
$ws.Range("L131").Value = 6690
$ws.Range("M131").Value = 781.5
$ws.Range("N131").Value = -16770
# Row 137: Creative Chocolate / Gateau au Chocolat
$ws.Range("H137").Value = 108414.31
$ws.Range("J137").Value = 95061.82000000001
$ws.Range("L137").Value = 285185.46
$ws.Range("N137").Value = -295385.46

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 20: Brothers in Arms / Brass Wristlets of Crafting
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
# Row 113: Copious Crystal Cannons / Manasilver Nugget
$ws.Range("H113").Value = 4851.6123
$ws.Range("I113").Value = 3554.919
$ws.Range("J113").Value = 8849.75
$ws.Range("K113").Value = 3554.919
$ws.Range("L113").Value = 8849.75
$ws.Range("M113").Value = -1384.919
$ws.Range("N113").Value = -13189.75
# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 3941.6843
$ws.Range("J132").Value = 4838
$ws.Range("L132").Value = 14514
$ws.Range("N132").Value = -19574

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 61: Spelling Me Softly / Raptor Leather
$ws.Range("H61").Value = 7550.76
$ws.Range("I61").Value = 7034.909
$ws.Range("K61").Value = 7034.909
$ws.Range("M61").Value = -6832.909
# Row 113: Peace in Rest / Atrociraptor Leather
$ws.Range("H113").Value = 7550.76
$ws.Range("I113").Value = 7034.909
$ws.Range("K113").Value = 7034.909
$ws.Range("M113").Value = -4864.909

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 22: Better Shroud than Sorry / Cotton Kurta
$ws.Range("H22").Value = 7000
$ws.Range("J22").Value = 6000
$ws.Range("L22").Value = 6000
$ws.Range("N22").Value = -6586
# Row 87: He Wears the Pants / Chimerical Felt Trousers
$ws.Range("H87").Value = 70000
$ws.Range("J87").Value = 70000
$ws.Range("L87").Value = 70000
$ws.Range("N87").Value = -72496
# Row 88: The Hat List / Chimerical Felt Cap of Scouting
$ws.Range("H88").Value = 70000
$ws.Range("J88").Value = 70000
$ws.Range("L88").Value = 70000
$ws.Range("N88").Value = -70812
# Row 90: Pom Hemlock (L) / Chimerical Felt Trousers
$ws.Range("H90").Value = 70000
$ws.Range("J90").Value = 70000
$ws.Range("L90").Value = 210000
$ws.Range("N90").Value = -222480
# Row 91: Knight Incognito (L) / Chimerical Felt Cap of Scouting
$ws.Range("H91").Value = 70000
$ws.Range("J91").Value = 70000
$ws.Range("L91").Value = 70000
$ws.Range("N91").Value = -72808
# Row 113: A Tender Table / Pixie Floss
$ws.Range("H113").Value = 1261.2122
$ws.Range("I113").Value = 1089.1578
$ws.Range("K113").Value = 3267.4734
$ws.Range("M113").Value = -1097.4734
# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 4122.5713
$ws.Range("J132").Value = 4391.2
$ws.Range("L132").Value = 13173.6
$ws.Range("N132").Value = -18233.6
# Row 133: Begin with the Basics / Snow Cotton Jacket
$ws.Range("H133").Value = 161000
$ws.Range("J133").Value = 161000
$ws.Range("L133").Value = 161000
$ws.Range("N133").Value = -171120
# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 21493100
$ws.Range("I136").Value = 38462744
$ws.Range("J136").Value = 483064.9
$ws.Range("K136").Value = 115388232
$ws.Range("L136").Value = 1449194.7
$ws.Range("M136").Value = -115385682
$ws.Range("N136").Value = -1454294.7
